# Reorders the EC (Estado de Cuenta) worker/period data block: the
# previously-second worker (EDGAR ELIECER CAMPO BORRE, 22 periods) now
# appears first (rows 16-37) and the previously-first worker (JONATHAN
# JOSE MUJICA PADILLA, 3 periods) now appears last (rows 38-40), with
# their "Valor Mora" / "Salario Basico" figures moving along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Force text storage for the document-number / period columns (C, E) so
# values like "73579467" or "2108" aren't auto-coerced to numbers. Column
# D (worker name) is left untouched since it already stores as General
# text and its style must not change.
$ws.Range("C16:C40").NumberFormat = "@"
$ws.Range("E16:E40").NumberFormat = "@"

$periodsEdgar = @("2108","2107","2106","2105","2104","2103","2102","2101", `
                   "2012","2011","2010","2009","2008","2007","2006","2005", `
                   "2004","2003","2002","2001","1912","1911")

$row = 16
foreach ($p in $periodsEdgar) {
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "73579467"
    $ws.Cells.Item($row, 4).Value = "EDGAR ELIECER CAMPO BORRE"
    $ws.Cells.Item($row, 5).Value = $p
    if ($p -eq "2108") {
        $ws.Cells.Item($row, 6).Value = 116000
    } else {
        $ws.Cells.Item($row, 6).Value = 120000
    }
    $ws.Cells.Item($row, 7).Value = 3000000
    $row++
}

$periodsJonathan = @("1908","1907","1906")
foreach ($p in $periodsJonathan) {
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "7920340"
    $ws.Cells.Item($row, 4).Value = "JONATHAN JOSE MUJICA PADILLA"
    $ws.Cells.Item($row, 5).Value = $p
    if ($p -eq "1906") {
        $ws.Cells.Item($row, 6).Value = 90000
    } else {
        $ws.Cells.Item($row, 6).Value = 100000
    }
    $ws.Cells.Item($row, 7).Value = 2500000
    $row++
}
